$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-11-09 Saturday" "2024-11-10 Sunday"

Replace-Text "520×2=" "349×6="
Replace-Text "628×6=" "410×6="
Replace-Text "155×8=" "986×9="
Replace-Text "980×3=" "919×7="
Replace-Text "712×2=" "868×4="

Replace-Text "287×2=" "346×6="
Replace-Text "963×7=" "947×5="
Replace-Text "529×7=" "945×4="
Replace-Text "626×6=" "412×2="
Replace-Text "688×7=" "547×7="

Replace-Text "727×9=" "869×5="
Replace-Text "772×5=" "820×6="
Replace-Text "431×6=" "315×5="
Replace-Text "282×4=" "335×8="
Replace-Text "425×6=" "501×6="

Replace-Text "805×7=" "112×3="
Replace-Text "818×7=" "695×2="
Replace-Text "249×2=" "263×6="
Replace-Text "869×8=" "814×7="
Replace-Text "486×5=" "410×5="

Replace-Text "853×6=" "917×6="
Replace-Text "878×4=" "444×3="
Replace-Text "532×3=" "568×4="
Replace-Text "310×9=" "396×7="
Replace-Text "862×8=" "980×4="
